$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed new shared strings in the exact authoring order via a scratch row ---
$scratchRow = 1000
$ws.Cells.Item($scratchRow, 1).Value = "Absolutely"
$ws.Cells.Item($scratchRow, 2).Value = "It is right"
$ws.Cells.Item($scratchRow, 3).Value = "Right"
$ws.Cells.Item($scratchRow, 4).Value = "It is true"
$ws.Cells.Item($scratchRow, 5).Value = "this is not true"
$ws.Cells.Item($scratchRow, 6).Value = "Not true"
$ws.Cells.Item($scratchRow, 7).Value = "it is accurate"
$ws.Cells.Item($scratchRow, 8).Value = "This is accurate"
$ws.Cells.Item($scratchRow, 9).Value = "Absolutely correct"
$ws.Cells.Item($scratchRow, 10).Value = "This is wrong"
$ws.Cells.Item($scratchRow, 11).Value = "It is wrong"
$ws.Cells.Item($scratchRow, 12).Value = "wrong"
$ws.Cells.Item($scratchRow, 13).Value = "This is false"
$ws.Cells.Item($scratchRow, 14).Value = "It is false"
$ws.Cells.Item($scratchRow, 15).Value = "It's false"
$ws.Cells.Item($scratchRow, 16).Value = "It is inaccurate"
$ws.Cells.Item($scratchRow, 17).Value = "Inaccurate"
$ws.Cells.Item($scratchRow, 18).Value = "this is inaccurate"
$ws.Cells.Item($scratchRow, 19).Value = "False "
$ws.Cells.Item($scratchRow, 20).Value = "True "
$ws.Cells.Item($scratchRow, 21).Value = "definitely"

# --- Write the final values for rows 17-63 (new reservation-acceptance phrases merged in) ---
$ws.Range("A17").Value = "definitely"
$ws.Range("B17").Value = "ACCEPT"
$ws.Range("A18").Value = "Absolutely correct"
$ws.Range("B18").Value = "ACCEPT"
$ws.Range("A19").Value = "That is correct"
$ws.Range("B19").Value = "ACCEPT"
$ws.Range("A20").Value = "Indeed"
$ws.Range("B20").Value = "ACCEPT"
$ws.Range("A21").Value = "True "
$ws.Range("B21").Value = "ACCEPT"
$ws.Range("A22").Value = "This is accurate"
$ws.Range("B22").Value = "ACCEPT"
$ws.Range("A23").Value = "it is accurate"
$ws.Range("B23").Value = "ACCEPT"
$ws.Range("A24").Value = "It is true"
$ws.Range("B24").Value = "ACCEPT"
$ws.Range("A25").Value = "Right"
$ws.Range("B25").Value = "ACCEPT"
$ws.Range("A26").Value = "It is right"
$ws.Range("B26").Value = "ACCEPT"
$ws.Range("A27").Value = "Correct"
$ws.Range("B27").Value = "ACCEPT"
$ws.Range("A28").Value = "It is correct"
$ws.Range("B28").Value = "ACCEPT"
$ws.Range("A29").Value = "Absolutely"
$ws.Range("B29").Value = "ACCEPT"
$ws.Range("A30").Value = "Indeed, It's for me"
$ws.Range("B30").Value = "ACCEPT"
$ws.Range("A31").Value = "No, It's not for me"
$ws.Range("B31").Value = "REJECT"
$ws.Range("A32").Value = "No, it is for"
$ws.Range("B32").Value = "REJECT"
$ws.Range("A33").Value = "Nope it's for"
$ws.Range("B33").Value = "REJECT"
$ws.Range("A34").Value = "it's for"
$ws.Range("B34").Value = "REJECT"
$ws.Range("A35").Value = "No"
$ws.Range("B35").Value = "REJECT"
$ws.Range("A36").Value = "No no it's for"
$ws.Range("B36").Value = "REJECT"
$ws.Range("A37").Value = "This is for"
$ws.Range("B37").Value = "REJECT"
$ws.Range("A38").Value = "it's not for me"
$ws.Range("B38").Value = "REJECT"
$ws.Range("A39").Value = "It is not for me. It's for"
$ws.Range("B39").Value = "REJECT"
$ws.Range("A40").Value = "This is not for me. It is for"
$ws.Range("B40").Value = "REJECT"
$ws.Range("A41").Value = "Nope"
$ws.Range("B41").Value = "REJECT"
$ws.Range("A42").Value = "No it wan't for me"
$ws.Range("B42").Value = "REJECT"
$ws.Range("A43").Value = "Not for me actually"
$ws.Range("B43").Value = "REJECT"
$ws.Range("A44").Value = "Well it's not for me"
$ws.Range("B44").Value = "REJECT"
$ws.Range("A45").Value = "No actually it was not for me"
$ws.Range("B45").Value = "REJECT"
$ws.Range("A46").Value = "It is not for me actually"
$ws.Range("B46").Value = "REJECT"
$ws.Range("A47").Value = "it is not correct"
$ws.Range("B47").Value = "REJECT"
$ws.Range("A48").Value = "Not Correct"
$ws.Range("B48").Value = "REJECT"
$ws.Range("A49").Value = "Incorrect"
$ws.Range("B49").Value = "REJECT"
$ws.Range("A50").Value = "That is not correct"
$ws.Range("B50").Value = "REJECT"
$ws.Range("A51").Value = "That is incorrect"
$ws.Range("B51").Value = "REJECT"
$ws.Range("A52").Value = "this is not true"
$ws.Range("B52").Value = "REJECT"
$ws.Range("A53").Value = "Not true"
$ws.Range("B53").Value = "REJECT"
$ws.Range("A54").Value = "This is wrong"
$ws.Range("B54").Value = "REJECT"
$ws.Range("A55").Value = "It is wrong"
$ws.Range("B55").Value = "REJECT"
$ws.Range("A56").Value = "wrong"
$ws.Range("B56").Value = "REJECT"
$ws.Range("A57").Value = "This is false"
$ws.Range("B57").Value = "REJECT"
$ws.Range("A58").Value = "False "
$ws.Range("B58").Value = "REJECT"
$ws.Range("A59").Value = "It is false"
$ws.Range("B59").Value = "REJECT"
$ws.Range("A60").Value = "It's false"
$ws.Range("B60").Value = "REJECT"
$ws.Range("A61").Value = "It is inaccurate"
$ws.Range("B61").Value = "REJECT"
$ws.Range("A62").Value = "Inaccurate"
$ws.Range("B62").Value = "REJECT"
$ws.Range("A63").Value = "this is inaccurate"
$ws.Range("B63").Value = "REJECT"

# --- Remove the scratch row used to seed shared-string order ---
$ws.Rows($scratchRow).Delete()

# --- Apply left-horizontal alignment to the cells the author left-aligned ---
$ws.Range("A21").HorizontalAlignment = -4131
$ws.Range("A22").HorizontalAlignment = -4131
$ws.Range("A23").HorizontalAlignment = -4131
$ws.Range("A58").HorizontalAlignment = -4131

# --- Restore the scroll position / active selection captured in the saved view ---
$ws.Range("A20").Select()
$excel.ActiveWindow.ScrollRow = 20
$ws.Range("F18").Select()